$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 1329
$ws.Range("C14").Value = "V1"

$ws.Range("A15").Value = 1329
$ws.Range("C15").Value = "LM"

$ws.Range("A16").Value = 1329
$ws.Range("C16").Value = "LI"

$ws.Range("C17").Select()
